$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.97"
$ws.Range("E2").Value = "'2.56%"
$ws.Range("D3").Value = "'35.63"
$ws.Range("E3").Value = "'12.49%"
$ws.Range("D4").Value = "'5.096"
$ws.Range("E4").Value = "'2.55%"
$ws.Range("D5").Value = "'0.07802"
$ws.Range("E5").Value = "'2.40%"
$ws.Range("D6").Value = "'2.267"
$ws.Range("E6").Value = "'1.33%"
$ws.Range("D7").Value = "'8.127"
$ws.Range("E7").Value = "'3.46%"
$ws.Range("D8").Value = "'4.029"
$ws.Range("E8").Value = "'6.55%"
$ws.Range("D9").Value = "'0.9275"
$ws.Range("E9").Value = "'0.00%"
$ws.Range("D10").Value = "'0.09658"
$ws.Range("E10").Value = "'-0.52%"
$ws.Range("D11").Value = "'0.1829"
$ws.Range("E11").Value = "'4.96%"
$ws.Range("D12").Value = "'0.08758"
$ws.Range("E12").Value = "'4.05%"
$ws.Range("D13").Value = "'0.03422"
$ws.Range("E13").Value = "'5.24%"
$ws.Range("D14").Value = "'0.09952"
$ws.Range("E14").Value = "'0.64%"
$ws.Range("D15").Value = "'0.001490"
$ws.Range("E15").Value = "'1.12%"
$ws.Range("D16").Value = "'0.005733"
$ws.Range("E16").Value = "'0.28%"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'-0.49%"
$ws.Range("E18").Value = "'-2.24%"
$ws.Range("D19").Value = "'0.3455"
$ws.Range("E19").Value = "'3.05%"
$ws.Range("D20").Value = "'0.1322"
$ws.Range("E20").Value = "'0.09%"
$ws.Range("D21").Value = "'4.586"
$ws.Range("E21").Value = "'12.56%"
$ws.Range("D22").Value = "'0.2237"
$ws.Range("E22").Value = "'-1.89%"
$ws.Range("D23").Value = "'0.04683"
$ws.Range("E23").Value = "'3.61%"
$ws.Range("E24").Value = "'2.69%"
$ws.Range("D25").Value = "'0.004541"
$ws.Range("E25").Value = "'4.58%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'0.84%"
$ws.Range("E27").Value = "'-19.93%"
$ws.Range("E39").Value = "'4.54%"
$ws.Range("D40").Value = "'0.04715"
$ws.Range("E40").Value = "'1.89%"
$ws.Range("D41").Value = "'0.007883"
$ws.Range("E41").Value = "'4.93%"
$ws.Range("D42").Value = "'0.1420"
$ws.Range("E42").Value = "'2.36%"
$ws.Range("D43").Value = "'0.008010"
$ws.Range("E43").Value = "'-17.96%"
$ws.Range("D44").Value = "'0.002294"
$ws.Range("E44").Value = "'7.02%"
$ws.Range("D45").Value = "'0.009111"
$ws.Range("E45").Value = "'-3.35%"
$ws.Range("D46").Value = "'0.00006235"
$ws.Range("E46").Value = "'2.97%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.83%"
$ws.Range("D48").Value = "'4.026"
$ws.Range("E48").Value = "'44.09%"
$ws.Range("D49").Value = "'0.002693"
$ws.Range("E49").Value = "'35.53%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.83%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.83%"
